# Auto-generated edit script applying the cryptos.xlsx price/volume update
# (GitHub Actions data refresh, commit Sun Sep  3 15:50:45 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '25.883.43'
$ws.Range('E2').Value = '  -0.39%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.637.14'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '214.47'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2549'
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06363'
$ws.Range('E9').Value = '  -1.02%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.44'
$ws.Range('E10').Value = '  -1.30%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07751'
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.653.61'
$ws.Range('E12').Value = '  +0.42%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.268'
$ws.Range('E13').Value = '  -0.22%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.5429'
$ws.Range('E14').Value = '  -0.73%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0₅7798'
$ws.Range('E15').Value = '  -2.01%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '64.11'
$ws.Range('E16').Value = '  -0.86%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '25.917.85'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '1.003'
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '196.04'
$ws.Range('E19').Value = '  -2.89%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.443'
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '9.903'
$ws.Range('E21').Value = '  -0.21%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.006'
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.006'
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.886'
$ws.Range('E24').Value = '  +0.35%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '140.99'
$ws.Range('E25').Value = '  -0.16%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1186'
$ws.Range('E26').Value = '  +3.98%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '6.843'
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.67'
$ws.Range('E28').Value = '  -0.47%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.234'
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.04924'
$ws.Range('E30').Value = '  -0.41%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.245'
$ws.Range('E31').Value = '  -1.01%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.166'
$ws.Range('E32').Value = '  -1.66%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.535'
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.365'
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.8918'
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '2.579'
$ws.Range('E36').Value = '  -1.85%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.132.33'
$ws.Range('E37').Value = '  -2.36%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.5410'
$ws.Range('E38').Value = '  -3.45%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01550'
$ws.Range('E39').Value = '  -1.11%  '
$ws.Range('E40').Value = '  -0.30%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.546'
$ws.Range('E41').Value = '  -0.45%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.571'
$ws.Range('E42').Value = '  -2.58%  '
$ws.Range('B43').Value = 'BabyDogeCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.0₈126'
$ws.Range('E43').Value = '  +9.51%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.8120'
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '99.34'
$ws.Range('E45').Value = '  -0.48%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.775.69'
$ws.Range('E46').Value = '  -0.35%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.4536'
$ws.Range('E47').Value = '  +0.17%  '
$ws.Range('E48').Value = '  -0.08%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '54.64'
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05065'
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.005'
$ws.Range('E51').Value = '  +0.06%  '
